# Generate Report for Handoff
# Updates the localization-status report: the item has moved from
# "In Translation" to "Ready for handoff", and the handoff timestamps
# are refreshed to the new generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn / de-de + generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-21 01:03:42"

# zh-cn sheet: Status + Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-21 01:03:38"

# de-de sheet: Status only (handoff datetime unchanged for this locale)
$dede.Range("C2").Value = "Ready for handoff"

# Columns widened to fit the new, longer status text ("Ready for handoff").
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
